# push testcases from 2 to 15
# Adds 4 new worksheets (Login, Incorrect_login, profile, Price_range_handler)
# with their test-fixture data, mirroring the author's manual Excel edit.

$wb = $excel.ActiveWorkbook

function Add-SheetAfterLast([string]$name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $last)
    $newSheet.Name = $name
    return $newSheet
}

# --- Login ---------------------------------------------------------------
$loginSheet = Add-SheetAfterLast "Login"
$loginSheet.Range("A1").Value = 9131899877
$loginSheet.Range("B1").Value = 9131899877
$loginSheet.Range("C1").Value = "zafar"
$loginSheet.Columns.Item(1).ColumnWidth = 33.21875
$loginSheet.Columns.Item(2).ColumnWidth = 30.21875
$loginSheet.Range("A1:B1").Select() | Out-Null

# --- Incorrect_login -------------------------------------------------------
$incorrectLoginSheet = Add-SheetAfterLast "Incorrect_login"
$incorrectLoginSheet.Range("A1").Value = 9131899877
$incorrectLoginSheet.Range("B1").Value = 9131899856
$incorrectLoginSheet.Range("C1").Font.Name = "Consolas"
$incorrectLoginSheet.Range("C1").Font.Family = 3
$incorrectLoginSheet.Range("C1").Font.Size = 10
$incorrectLoginSheet.Range("C1").Font.Color = 16711722
$incorrectLoginSheet.Columns.Item(1).ColumnWidth = 23.5546875
$incorrectLoginSheet.Columns.Item(2).ColumnWidth = 38.5546875
$incorrectLoginSheet.Columns.Item(3).ColumnWidth = 38.5546875
$incorrectLoginSheet.Columns.Item(4).ColumnWidth = 28.5546875
$incorrectLoginSheet.PageSetup.Orientation = 1
$incorrectLoginSheet.Range("D1").Select() | Out-Null

# --- profile ---------------------------------------------------------------
$profileSheet = Add-SheetAfterLast "profile"
$profileSheet.Range("A1").Value = "zafar shareef"
$profileSheet.Columns.Item(1).ColumnWidth = 15.6640625
$profileSheet.Range("A2").Select() | Out-Null

# --- Price_range_handler ----------------------------------------------------
$priceRangeSheet = Add-SheetAfterLast "Price_range_handler"
$priceRangeSheet.Range("A1").Value = 35
$priceRangeSheet.Range("B1").Value = "Sport & Outdoor"
$priceRangeSheet.Columns.Item(2).ColumnWidth = 19.33203125
$priceRangeSheet.Range("C3").Select() | Out-Null

Write-Output "done"
